$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values per the commit diff (cryptos list refresh)
# Cells whose new text looks like a plain number must be forced to stay
# text (matching the original inlineStr cells) instead of being coerced
# to a numeric type by the COM layer's automatic type inference.

$ws.Range('D2').Value = '71.151.35'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '2.624.93'
$ws.Range('E3').Value = '  +4.49%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.36%  '
$ws.Range('D9').Value = '2.624.98'
$ws.Range('E10').Value = '  +15.31%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '3.074.29'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('E15').Value = '  +3.04%  '
$ws.Range('E16').Value = '  +7.54%  '
$ws.Range('D17').Value = '71.120.66'
$ws.Range('E17').Value = '  +4.05%  '
$ws.Range('D18').Value = '2.627.57'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '382.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  +5.76%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +8.26%  '
$ws.Range('E27').Value = '  +4.06%  '
$ws.Range('D28').Value = '2.759.70'
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '0.0₃0956'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '539.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('E33').Value = '  +4.36%  '
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('E36').Value = '  +2.31%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.21'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.56%  '
$ws.Range('E39').Value = '  +7.59%  '
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  +4.10%  '
$ws.Range('E42').Value = '  +9.23%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.331'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '154.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.533'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.11%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0265'
$ws.Range('E51').Value = '  +1.18%  '
